$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '49.552.56' },
    @{ Cell = 'E2'; Value = '  -0.66%  ' },
    @{ Cell = 'D3'; Value = '2.640.35' },
    @{ Cell = 'E3'; Value = '  -0.35%  ' },
    @{ Cell = 'D4'; Value = '0.999' },
    @{ Cell = 'E4'; Value = '  -0.04%  ' },
    @{ Cell = 'D5'; Value = '111.69' },
    @{ Cell = 'E5'; Value = '  -1.77%  ' },
    @{ Cell = 'D6'; Value = '325.66' },
    @{ Cell = 'E6'; Value = '  -0.53%  ' },
    @{ Cell = 'D7'; Value = '0.524' },
    @{ Cell = 'E7'; Value = '  -1.12%  ' },
    @{ Cell = 'D8'; Value = '0.999' },
    @{ Cell = 'E8'; Value = '  -0.02%  ' },
    @{ Cell = 'D9'; Value = '0.546' },
    @{ Cell = 'E9'; Value = '  -1.43%  ' },
    @{ Cell = 'D10'; Value = '39.53' },
    @{ Cell = 'E10'; Value = '  -4.02%  ' },
    @{ Cell = 'D11'; Value = '20.10' },
    @{ Cell = 'E11'; Value = '  -0.52%  ' },
    @{ Cell = 'D12'; Value = '0.0810' },
    @{ Cell = 'E12'; Value = '  -1.49%  ' },
    @{ Cell = 'E13'; Value = '  +1.48%  ' },
    @{ Cell = 'D14'; Value = '7.54' },
    @{ Cell = 'E14'; Value = '  +2.67%  ' },
    @{ Cell = 'D15'; Value = '3.049.85' },
    @{ Cell = 'E15'; Value = '  -0.66%  ' },
    @{ Cell = 'D16'; Value = '2.637.98' },
    @{ Cell = 'E16'; Value = '  -0.47%  ' },
    @{ Cell = 'D17'; Value = '0.853' },
    @{ Cell = 'E17'; Value = '  -2.20%  ' },
    @{ Cell = 'D18'; Value = '49.491.92' },
    @{ Cell = 'E18'; Value = '  -0.68%  ' },
    @{ Cell = 'D19'; Value = '13.09' },
    @{ Cell = 'E19'; Value = '  -0.76%  ' },
    @{ Cell = 'D20'; Value = '2.92' },
    @{ Cell = 'E20'; Value = '  -0.54%  ' },
    @{ Cell = 'E21'; Value = '  -1.44%  ' },
    @{ Cell = 'D22'; Value = '0.0₃0948' },
    @{ Cell = 'E22'; Value = '  -0.92%  ' },
    @{ Cell = 'D23'; Value = '268.70' },
    @{ Cell = 'E23'; Value = '  -2.90%  ' },
    @{ Cell = 'D24'; Value = '69.11' },
    @{ Cell = 'E24'; Value = '  -4.33%  ' },
    @{ Cell = 'D25'; Value = '2.55' },
    @{ Cell = 'E25'; Value = '  -1.50%  ' },
    @{ Cell = 'B26'; Value = 'EthereumClassic' },
    @{ Cell = 'C26'; Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc' },
    @{ Cell = 'D26'; Value = '26.09' },
    @{ Cell = 'E26'; Value = '  -2.63%  ' },
    @{ Cell = 'B27'; Value = 'Dai' },
    @{ Cell = 'C27'; Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai' },
    @{ Cell = 'D27'; Value = '1.00' },
    @{ Cell = 'E27'; Value = '  +0.07%  ' },
    @{ Cell = 'D28'; Value = '10.18' },
    @{ Cell = 'E28'; Value = '  +1.74%  ' },
    @{ Cell = 'E29'; Value = '  -1.11%  ' },
    @{ Cell = 'D30'; Value = '0.138' },
    @{ Cell = 'E30'; Value = '  -1.18%  ' },
    @{ Cell = 'D31'; Value = '34.69' },
    @{ Cell = 'E31'; Value = '  -3.61%  ' },
    @{ Cell = 'D32'; Value = '49.57' },
    @{ Cell = 'E32'; Value = '  -1.53%  ' },
    @{ Cell = 'D33'; Value = '5.49' },
    @{ Cell = 'E33'; Value = '  +1.35%  ' },
    @{ Cell = 'D34'; Value = '0.0807' },
    @{ Cell = 'E34'; Value = '  -0.09%  ' },
    @{ Cell = 'E35'; Value = '  -0.19%  ' },
    @{ Cell = 'D36'; Value = '18.99' },
    @{ Cell = 'E36'; Value = '  -3.19%  ' },
    @{ Cell = 'D37'; Value = '4.94' },
    @{ Cell = 'E37'; Value = '  +3.50%  ' },
    @{ Cell = 'D38'; Value = '2.03' },
    @{ Cell = 'E38'; Value = '  -2.18%  ' },
    @{ Cell = 'D39'; Value = '3.09' },
    @{ Cell = 'E39'; Value = '  +0.75%  ' },
    @{ Cell = 'D40'; Value = '129.14' },
    @{ Cell = 'E40'; Value = '  +3.06%  ' },
    @{ Cell = 'D41'; Value = '22.67' },
    @{ Cell = 'E41'; Value = '  +1.74%  ' },
    @{ Cell = 'E42'; Value = '  -1.27%  ' },
    @{ Cell = 'D43'; Value = '2.23' },
    @{ Cell = 'E43'; Value = '  +0.05%  ' },
    @{ Cell = 'D44'; Value = '0.0329' },
    @{ Cell = 'E44'; Value = '  +4.43%  ' },
    @{ Cell = 'D45'; Value = '2.047.37' },
    @{ Cell = 'E45'; Value = '  -1.36%  ' },
    @{ Cell = 'E46'; Value = '  +8.65%  ' },
    @{ Cell = 'D47'; Value = '3.26' },
    @{ Cell = 'E47'; Value = '  -2.40%  ' },
    @{ Cell = 'E48'; Value = '  -4.66%  ' },
    @{ Cell = 'D49'; Value = '8.82' },
    @{ Cell = 'E49'; Value = '  -3.40%  ' },
    @{ Cell = 'D50'; Value = '5.22' },
    @{ Cell = 'E50'; Value = '  -3.37%  ' },
    @{ Cell = 'D51'; Value = '58.70' },
    @{ Cell = 'E51'; Value = '  +0.08%  ' }
)

foreach ($u in $updates) {
    $c = $ws.Range($u.Cell)
    $c.Value = ("'" + $u.Value)
    $c.Style = 'Normal'
}
